$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030600518318875
$ws.Range("D2").Value = 1.03709337814674
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.04448852808006
$ws.Range("I2").Value = 1.032170972291224
$ws.Range("J2").Value = 1.035740866140009
$ws.Range("K2").Value = 1.039885246867717
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.047259436336952
$ws.Range("N2").Value = 1.015849709895317

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031743181901406
$ws.Range("D3").Value = 1.037951745522706
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.04555536072421
$ws.Range("I3").Value = 1.032357850114595
$ws.Range("J3").Value = 1.036524066865788
$ws.Range("K3").Value = 1.040553222017691
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.048136850418953
$ws.Range("N3").Value = 1.016113922389645

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032482367351502
$ws.Range("D4").Value = 1.038506845874865
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.046245789504864
$ws.Range("I4").Value = 1.032477263897701
$ws.Range("J4").Value = 1.037030140089224
$ws.Range("K4").Value = 1.040984487503197
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.048704116825061
$ws.Range("N4").Value = 1.016284520699327

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032793075585139
$ws.Range("D5").Value = 1.038740133289622
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.046536074035559
$ws.Range("I5").Value = 1.03252710422103
$ws.Range("J5").Value = 1.037242723728576
$ws.Range("K5").Value = 1.041165562063588
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.048942481188203
$ws.Range("N5").Value = 1.016356152953226

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032845242241062
$ws.Range("D6").Value = 1.038779298771329
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.046584815752911
$ws.Range("I6").Value = 1.032535451446281
$ws.Range("J6").Value = 1.037278407527539
$ws.Range("K6").Value = 1.041195951828642
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.048982496911312
$ws.Range("N6").Value = 1.016368175207937

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032486519226758
$ws.Range("D7").Value = 1.038509963374382
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.04624966818948
$ws.Range("I7").Value = 1.032477931286002
$ws.Range("J7").Value = 1.037032981307871
$ws.Range("K7").Value = 1.040986907931725
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.048707302310652
$ws.Range("N7").Value = 1.016285478195859

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030986727770101
$ws.Range("D8").Value = 1.03738353383967
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.044849045196679
$ws.Range("I8").Value = 1.032234440873709
$ws.Range("J8").Value = 1.036005699769019
$ws.Range("K8").Value = 1.040111190883766
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.04755606192077
$ws.Range("N8").Value = 1.015939077329283

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028342362038526
$ws.Range("D9").Value = 1.035396163278834
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.04238184709577
$ws.Range("I9").Value = 1.031793827151251
$ws.Range("J9").Value = 1.034190040475406
$ws.Range("K9").Value = 1.038560712108576
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.045523746051347
$ws.Range("N9").Value = 1.015325877169126

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026578337591317
$ws.Range("D10").Value = 1.034069591982118
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.040737614229483
$ws.Range("I10").Value = 1.031492318917747
$ws.Range("J10").Value = 1.032975900611369
$ws.Range("K10").Value = 1.03752210188491
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.044166370827932
$ws.Range("N10").Value = 1.014915190415995

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025814213214334
$ws.Range("D11").Value = 1.03349477665957
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.040025770902133
$ws.Range("I11").Value = 1.031359920475243
$ws.Range("J11").Value = 1.032449279317296
$ws.Range("K11").Value = 1.037071192408901
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.043578014134873
$ws.Range("N11").Value = 1.014736909423104

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025530338140027
$ws.Range("D12").Value = 1.033281204080355
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.039761377995992
$ws.Range("I12").Value = 1.031310464774749
$ws.Range("J12").Value = 1.032253534138298
$ws.Range("K12").Value = 1.036903526172632
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.043359380871759
$ws.Range("N12").Value = 1.014670619999527

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02559123238192
$ws.Range("D13").Value = 1.033327018869593
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.039818090395449
$ws.Range("I13").Value = 1.031321085727533
$ws.Range("J13").Value = 1.032295528246088
$ws.Range("K13").Value = 1.036939499228627
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.043406282594283
$ws.Range("N13").Value = 1.014684842382555

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025790748945035
$ws.Range("D14").Value = 1.033477123921007
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.040003915751026
$ws.Range("I14").Value = 1.031355838103258
$ws.Range("J14").Value = 1.0324331017096
$ws.Range("K14").Value = 1.037057336706956
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.04355994370901
$ws.Range("N14").Value = 1.014731431305364

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025913671637422
$ws.Range("D15").Value = 1.033569600561477
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.040118411170878
$ws.Range("I15").Value = 1.031377213478272
$ws.Range("J15").Value = 1.032517847404617
$ws.Range("K15").Value = 1.037129916616206
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.043654607264862
$ws.Range("N15").Value = 1.01476012727223

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026629043754353
$ws.Range("D16").Value = 1.03410773210378
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.040784859404392
$ws.Range("I16").Value = 1.031501066921838
$ws.Range("J16").Value = 1.033010831900875
$ws.Range("K16").Value = 1.037552002270653
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.044205405341099
$ws.Range("N16").Value = 1.014927012822733

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027077698927271
$ws.Range("D17").Value = 1.034445180181926
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.041202936517049
$ws.Range("I17").Value = 1.031578263334675
$ws.Range("J17").Value = 1.033319829043498
$ws.Range("K17").Value = 1.037816448046458
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.04455074421717
$ws.Range("N17").Value = 1.015031574869881

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027339364015654
$ws.Range("D18").Value = 1.034641969247473
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.041446805481976
$ws.Range("I18").Value = 1.031623112841199
$ws.Range("J18").Value = 1.033499975958671
$ws.Range("K18").Value = 1.037970580546077
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.044752116279176
$ws.Range("N18").Value = 1.015092520647981

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027428580369887
$ws.Range("D19").Value = 1.034709062679986
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.041529960494995
$ws.Range("I19").Value = 1.031638375198892
$ws.Range("J19").Value = 1.033561386855956
$ws.Range("K19").Value = 1.038023116380236
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.044820769084588
$ws.Range("N19").Value = 1.015113294211259

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027029565369134
$ws.Range("D20").Value = 1.034408979200537
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.041158079598349
$ws.Range("I20").Value = 1.031569999287024
$ws.Range("J20").Value = 1.033286685463781
$ws.Range("K20").Value = 1.03778808733756
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.044513698654516
$ws.Range("N20").Value = 1.015020360851912

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02573199755276
$ws.Range("D21").Value = 1.033432923399229
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.039949194340543
$ws.Range("I21").Value = 1.031345612042619
$ws.Range("J21").Value = 1.032392593474801
$ws.Range("K21").Value = 1.037022641434316
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.043514696868354
$ws.Range("N21").Value = 1.014717713905305

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02491590248975
$ws.Range("D22").Value = 1.032818887758056
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.03918922053135
$ws.Range("I22").Value = 1.031202928062783
$ws.Range("J22").Value = 1.031829663321579
$ws.Range("K22").Value = 1.036540342941994
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.042886055732805
$ws.Range("N22").Value = 1.014527034675364

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025348555379861
$ws.Range("D23").Value = 1.033144433005657
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.039592087745303
$ws.Range("I23").Value = 1.031278719494047
$ws.Range("J23").Value = 1.032128157345748
$ws.Range("K23").Value = 1.036796116493217
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.043219360646014
$ws.Range("N23").Value = 1.014628154663105

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027051314937846
$ws.Range("D24").Value = 1.034425336985115
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.041178348472756
$ws.Range("I24").Value = 1.031573734003899
$ws.Range("J24").Value = 1.03330166188591
$ws.Range("K24").Value = 1.037800902673928
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.044530438129427
$ws.Range("N24").Value = 1.01502542811824

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029026184616443
$ws.Range("D25").Value = 1.035910237969554
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.043019575366228
$ws.Range("I25").Value = 1.031909105229971
$ws.Range("J25").Value = 1.034660081866403
$ws.Range("K25").Value = 1.038962420846856
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.046049587102269
$ws.Range("N25").Value = 1.015484736215876
